# Lightning SDK1 To-Do schedule update
# - begin working on debug, implement custom services, fix ServiceNotificationType.Shutdown not actually working

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Rename "OnTrigger animations" to "OnEvent / FromScript animations" (row 18) ---
$ws.Cells.Item(18, 1).Value() = "OnEvent / FromScript animations"

# --- 2. Mark OnObjectCreated / OnObjectDestroyed events complete (rows 8 & 9) ---
$cComplete = $ws.Cells.Item(24, 3)
$dDate = $ws.Cells.Item(24, 4)

$cComplete.Copy()
$ws.Cells.Item(8, 3).PasteSpecial(-4122)
$ws.Cells.Item(8, 3).Value() = "Complete"
$dDate.Copy()
$ws.Cells.Item(8, 4).PasteSpecial(-4122)
$ws.Cells.Item(8, 4).Value() = "Aug 20, 2021"

$cComplete.Copy()
$ws.Cells.Item(9, 3).PasteSpecial(-4122)
$ws.Cells.Item(9, 3).Value() = "Complete"
$dDate.Copy()
$ws.Cells.Item(9, 4).PasteSpecial(-4122)
$ws.Cells.Item(9, 4).Value() = "Aug 20, 2021"

# --- 3. Mark Instance.GetInheritableClasses() complete (row 25) ---
$cComplete.Copy()
$ws.Cells.Item(25, 3).PasteSpecial(-4122)
$ws.Cells.Item(25, 3).Value() = "Complete"
$dDate.Copy()
$ws.Cells.Item(25, 4).PasteSpecial(-4122)
$ws.Cells.Item(25, 4).Value() = "Aug 19, 2021"

# --- 4. Insert a new row at 27 for "DDMS Objects know when fully loaded from XML" ---
$ws.Rows(27).Insert(-4121)
$ws.Cells.Item(27, 1).Value() = "DDMS Objects know when fully loaded from XML"
$ws.Cells.Item(27, 2).Value() = "Improvement"

# --- 5. Append two new "External - Documentation" rows (44 & 45) ---
$srcA = $ws.Cells.Item(4, 1)
$srcA.Copy()
$ws.Cells.Item(44, 1).PasteSpecial(-4122)
$ws.Cells.Item(44, 1).Value() = "Welcome Document"

$ws.Cells.Item(44, 2).Value() = "External - Documentation"
$ws.Cells.Item(44, 2).Font.Name = "Arial"
$ws.Cells.Item(44, 2).Font.Bold = $false
$ws.Cells.Item(44, 2).Interior.Color = 10498160

$srcA.Copy()
$ws.Cells.Item(45, 1).PasteSpecial(-4122)
$ws.Cells.Item(45, 1).Value() = "Lua Documentation"

$ws.Cells.Item(44, 2).Copy()
$ws.Cells.Item(45, 2).PasteSpecial(-4122)
$ws.Cells.Item(45, 2).Value() = "External - Documentation"

# --- 6. Rename Win32 dialog fix feature (now row 30) ---
$ws.Cells.Item(30, 1).Value() = "Fix Win32 Dialog File Names (Lightning.Core.NativeInterop.Win32)"

# --- 7. Update view state to match where the author ended up editing ---
$ws.Range("B33").Select()
